{"js": "// Fix the typo \"Habilidades t\u00e9nicas\" -> \"Habilidades t\u00e9cnicas\"\n// (missing \"c\" in \"t\u00e9cnicas\") in the \"Skills\" section heading.\n//\n// The original run's text \"Habilidades t\u00e9nicas\" is replaced with the\n// corrected \"Habilidades t\u00e9cnicas\" while keeping the same paragraph /\n// run formatting (style \"Heading 2\", <w:noProof/>, <w:lang w:bidi=\"es-ES\"/>).\n\nconst body = context.document.body;\n\n// Locate the exact (unique) misspelled heading text in the document body.\nconst matches = body.search(\"Habilidades t\u00e9nicas\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items,text\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"Habilidades t\u00e9nicas\" in the document.');\n}\n\n// Replace the whole matched range with the corrected text; this keeps the\n// run's existing character formatting (rPr) intact.\nfor (const range of matches.items) {\n  range.insertText(\"Habilidades t\u00e9cnicas\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix the typo \"Habilidades t\u00e9nicas\" -> \"Habilidades t\u00e9cnicas\"\n# (missing \"c\" in \"t\u00e9cnicas\") in the \"Skills\" section heading.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Habilidades t\u00e9nicas\"\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Habilidades t\u00e9cnicas\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw 'Could not find \"Habilidades t\u00e9nicas\" in the document.'\n}\n\nWrite-Output \"done\"\n"}
